$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.918.24"
$ws.Range("D2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.649.25"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.41%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.28%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.42%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.42%  "

$ws.Range("E7").Value = "  -1.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3832"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.64%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "51.64"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.72%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.349"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.41%  "

$ws.Range("E11").Value = "  +0.54%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08434"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.74%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.87"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.62%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.070"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.59%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.906"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.35%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001313"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.80%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.650.27"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.56%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.41"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.31%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06982"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.84%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.79%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.928"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.67%  "

$ws.Range("E22").Value = "  +0.50%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.69"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.35%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.878.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.47%  "

$ws.Range("E25").Value = "  -0.68%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.958"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.87%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.35%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "151.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.82%  "

$ws.Range("E29").Value = "  +1.56%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "138.30"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.77%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.829"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.19%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.519"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.09%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.829.01"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.59%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.049"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.84%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08028"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.19%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02944"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.79%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.692"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.28%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "10.89"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.98%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2672"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.44%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09093"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.92%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7547"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.78%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.45"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.98%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.417"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.20%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.34"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.69%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6952"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.34%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.455"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.23%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.079"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.14%  "

$ws.Range("E48").Value = "  +0.40%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08278"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.36%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "134.09"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.90%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.208"
$ws.Range("D51").Style = "Normal"

